# Update the cryptos list with freshly scraped price/volume data.
# Mirrors the GitHub Actions scraper commit that refreshed D (Price) and
# E (Volume(1h)) columns, plus a reordering of the Polkadot / WrappedEther
# rows (12 & 13) including their Coin name and Link columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 12 & 13 swap which coin (and its link) occupies which row, in
# addition to getting refreshed price/volume figures.
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"

# row => (Price, Volume(1h))
$updates = @{
    2  = @("25.879.06", "  +0.41%  ")
    3  = @("1.637.43",  "  -0.17%  ")
    4  = @($null,        "  -1.90%  ")
    5  = @("214.28",    "  -1.25%  ")
    6  = @("0.5036",    "  -0.55%  ")
    7  = @("1.001",     "  -1.86%  ")
    8  = @("0.2567",    "  -0.72%  ")
    9  = @("0.06371",   "  -0.80%  ")
    10 = @("19.45",     "  -0.38%  ")
    11 = @("0.07770",   "  -0.18%  ")
    12 = @("1.649.20",  "  -0.09%  ")
    13 = @("4.253",     "  -0.36%  ")
    14 = @("1.866.94",  "  +0.19%  ")
    15 = @("0.5412",    "  -1.16%  ")
    16 = @("0.0₅7893",  "  -0.65%  ")
    17 = @("64.36",     "  +1.20%  ")
    18 = @("25.923.28", "  -0.02%  ")
    19 = @("1.002",     "  -1.81%  ")
    20 = @("197.27",    "  -3.83%  ")
    21 = @("4.353",     "  -0.09%  ")
    22 = @("9.864",     "  -1.53%  ")
    23 = @("5.953",     "  -0.23%  ")
    24 = @("1.003",     "  -1.98%  ")
    25 = @("1.874",     "  -5.23%  ")
    26 = @("140.93",    "  -0.92%  ")
    27 = @("0.1133",    "  -2.00%  ")
    28 = @("6.816",     "  -0.15%  ")
    29 = @("15.66",     "  -0.59%  ")
    30 = @("1.235",     "  -1.14%  ")
    31 = @("0.04927",   "  -2.21%  ")
    32 = @("3.256",     "  -0.47%  ")
    33 = @("3.183",     "  -1.10%  ")
    34 = @("1.531",     "  -1.07%  ")
    35 = @("2.363",     "  -0.20%  ")
    36 = @("0.8890",    "  -0.90%  ")
    37 = @("2.603",     "  -2.22%  ")
    38 = @("1.142.98",  "  +1.64%  ")
    39 = @("0.5532",    "  -2.51%  ")
    40 = @("0.01565",   "  -0.03%  ")
    41 = @("1.001",     "  -2.42%  ")
    42 = @("5.684",     "  +0.23%  ")
    43 = @("0.8094",    "  -0.77%  ")
    44 = @("99.58",     "  -0.36%  ")
    45 = @($null,        "  +6.39%  ")
    46 = @("1.778.15",  "  +0.37%  ")
    47 = @("0.4509",    "  -1.02%  ")
    48 = @("0.9988",    "  -1.88%  ")
    49 = @("54.49",     "  -1.03%  ")
    50 = @("0.05056",   "  +0.12%  ")
    51 = @("1.005",     "  -1.05%  ")
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $price = $pair[0]
    $volume = $pair[1]

    if ($null -ne $price) {
        # Prefix with an apostrophe so Excel stores the digit-and-dot-grouped
        # price (e.g. "0.07770", "25.879.06") as text instead of silently
        # re-parsing/normalising it as a number and dropping trailing zeros.
        $ws.Cells.Item($row, 4).Value = "'" + $price
    }
    $ws.Cells.Item($row, 5).Value = $volume
}

Write-Host "Updated cryptos list"
